$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Etape 4 : checkClass() ---
# Values are entered in this specific order so that new shared strings
# land at the same indices as in the target workbook.
$ws.Range("A19").Value = "4 : checkClass()"

$ws.Range("E21").Value = "message erreur + nom + classe fautive (classe anormale et classe attendue) + stop"
$ws.Range("B21").Value = "classe incorrecte factor"
$ws.Range("C21").Value = "x"

$ws.Range("B22").Value = "classe incorrecte integer"
$ws.Range("C22").Value = "x"
$ws.Range("E22").Value = "message erreur + nom + classe fautive (classe anormale et classe attendue) + stop"

$ws.Range("B23").Value = "classe incorrecte numeric"
$ws.Range("C23").Value = "x"
$ws.Range("E23").Value = "message erreur + nom + classe fautive (classe anormale et classe attendue) + stop"

$ws.Range("B19").Value = "classe attendue IND"
$ws.Range("C19").Value = "x"

$ws.Range("B20").Value = "classe attendue QUAD"
$ws.Range("C20").Value = "x"

$ws.Range("E19").Value = "factor, integer ou numeric : message info"
$ws.Range("E20").Value = "factor, integer ou numeric : message info"

# --- Etape 5 : checkFactor() ---
$ws.Range("A24").Value = "5: checkFactor()"
$ws.Range("A24").NumberFormat = "h:mm"

$ws.Range("B24").Value = "uniquement les colonnes de classe factor"
$ws.Range("B25").Value = "modalité correspond à l'attendu"
$ws.Range("B26").Value = "modalité ne correspond pas à l'attendu"

$ws.Range("E26").Value = "message erreur + col/numero/contenu des lignes + contenu attendu"
$ws.Range("E25").Value = "message d'information pertinent"

# Restore the selected cell as it ended up in the saved workbook
$ws.Range("B15").Select()
